$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column C
$ws.Range("C1").Value = "function "

# Set width for new column C (Excel ColumnWidth units differ from the
# raw OOXML "width" attribute by an offset of ~0.83; 37.1667 yields width="38")
$ws.Columns.Item(3).ColumnWidth = 37.1667

# New rows of issue log entries
$ws.Range("C5").Value = "check "
$ws.Range("A5").Value = "some cards stay flipped even if they don't match "

$ws.Range("A6").Value = "flippping animation happens faster when clicking faster "

$ws.Range("A7").Value = "seems to match incorrectly when clicking faster "

$ws.Range("A8").Value = "flipping doesn't happen in the right order "

$ws.Range("A9").Value = "footer "
$ws.Range("B9").Value = "make it absolute position "

$ws.Range("A10").Value = "media queries "

$ws.Range("A11").Value = "stars function for the win page does not work properly "

# Update selection to match target state
$ws.Range("A11").Select()
